$wb = $excel.ActiveWorkbook

# --- Helper: write a value as TEXT (preserve leading zeros / trailing decimal zeros) ---
function Set-TextValue($cell, $val) {
    $cell.Value = "'" + $val
    $cell.ClearFormats()
}

# ============================================================
# Step 1: insert "2022-Q1" worksheet before "总计", fund-holding table
# ============================================================
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($lastSheet)
$ws.Name = '2022-Q1'

# Header row (copy the bold/centered/bordered style used by sibling sheets)
$srcHeader = $wb.Worksheets.Item('2021-Q4').Range("B1:H1")
$srcHeader.Copy()
$ws.Range("B1:H1").PasteSpecial(-4122)
$ws.Range("B1").Value = '基金代码'
$ws.Range("C1").Value = '基金名称'
$ws.Range("D1").Value = '基金规模'
$ws.Range("E1").Value = '股票总仓位'
$ws.Range("F1").Value = '仓位占比'
$ws.Range("G1").Value = '持有市值(亿元)'
$ws.Range("H1").Value = '仓位排名'

# Row-number style source (A column, bold/centered/bordered like header)
$aStyleSrc = $wb.Worksheets.Item('2021-Q4').Range("A2")

# Data rows
# row 2
$aStyleSrc.Copy()
$ws.Cells.Item(2,1).PasteSpecial(-4122)
$ws.Cells.Item(2,1).Value = 0
Set-TextValue $ws.Cells.Item(2,2) '900010'
Set-TextValue $ws.Cells.Item(2,3) '中信卓越成长两年持有期混合A'
Set-TextValue $ws.Cells.Item(2,4) '133.02'
Set-TextValue $ws.Cells.Item(2,5) '93.07'
Set-TextValue $ws.Cells.Item(2,6) '3.16'
Set-TextValue $ws.Cells.Item(2,7) '4.2034'
$ws.Cells.Item(2,8).Value = 7

# row 3
$aStyleSrc.Copy()
$ws.Cells.Item(3,1).PasteSpecial(-4122)
$ws.Cells.Item(3,1).Value = 1
Set-TextValue $ws.Cells.Item(3,2) '900090'
Set-TextValue $ws.Cells.Item(3,3) '中信卓越成长两年持有期混合B'
Set-TextValue $ws.Cells.Item(3,4) '86.95'
Set-TextValue $ws.Cells.Item(3,5) '93.07'
Set-TextValue $ws.Cells.Item(3,6) '3.16'
Set-TextValue $ws.Cells.Item(3,7) '2.7476'
$ws.Cells.Item(3,8).Value = 7

# row 4
$aStyleSrc.Copy()
$ws.Cells.Item(4,1).PasteSpecial(-4122)
$ws.Cells.Item(4,1).Value = 2
Set-TextValue $ws.Cells.Item(4,2) '012344'
Set-TextValue $ws.Cells.Item(4,3) '嘉实领先优势混合型证券投资基金A'
Set-TextValue $ws.Cells.Item(4,4) '68.76'
Set-TextValue $ws.Cells.Item(4,5) '82.99'
Set-TextValue $ws.Cells.Item(4,6) '3.01'
Set-TextValue $ws.Cells.Item(4,7) '2.0697'
$ws.Cells.Item(4,8).Value = 9

# row 5
$aStyleSrc.Copy()
$ws.Cells.Item(5,1).PasteSpecial(-4122)
$ws.Cells.Item(5,1).Value = 3
Set-TextValue $ws.Cells.Item(5,2) '169103'
Set-TextValue $ws.Cells.Item(5,3) '东方红睿轩三年定期开放灵活配置混合'
Set-TextValue $ws.Cells.Item(5,4) '40.03'
Set-TextValue $ws.Cells.Item(5,5) '57.34'
Set-TextValue $ws.Cells.Item(5,6) '3.16'
Set-TextValue $ws.Cells.Item(5,7) '1.2649'
$ws.Cells.Item(5,8).Value = 4

# row 6
$aStyleSrc.Copy()
$ws.Cells.Item(6,1).PasteSpecial(-4122)
$ws.Cells.Item(6,1).Value = 4
Set-TextValue $ws.Cells.Item(6,2) '004278'
Set-TextValue $ws.Cells.Item(6,3) '东方红智逸沪港深定期开放混合'
Set-TextValue $ws.Cells.Item(6,4) '34.95'
Set-TextValue $ws.Cells.Item(6,5) '22.66'
Set-TextValue $ws.Cells.Item(6,6) '1.98'
Set-TextValue $ws.Cells.Item(6,7) '0.6920'
$ws.Cells.Item(6,8).Value = 4

# row 7
$aStyleSrc.Copy()
$ws.Cells.Item(7,1).PasteSpecial(-4122)
$ws.Cells.Item(7,1).Value = 5
Set-TextValue $ws.Cells.Item(7,2) '003396'
Set-TextValue $ws.Cells.Item(7,3) '东方红优享红利沪港深灵活配置混合'
Set-TextValue $ws.Cells.Item(7,4) '15.87'
Set-TextValue $ws.Cells.Item(7,5) '89.27'
Set-TextValue $ws.Cells.Item(7,6) '4.29'
Set-TextValue $ws.Cells.Item(7,7) '0.6808'
$ws.Cells.Item(7,8).Value = 6

# row 8
$aStyleSrc.Copy()
$ws.Cells.Item(8,1).PasteSpecial(-4122)
$ws.Cells.Item(8,1).Value = 6
Set-TextValue $ws.Cells.Item(8,2) '160726'
Set-TextValue $ws.Cells.Item(8,3) '嘉实瑞享定期开放灵活配置混合'
Set-TextValue $ws.Cells.Item(8,4) '23.58'
Set-TextValue $ws.Cells.Item(8,5) '63.95'
Set-TextValue $ws.Cells.Item(8,6) '2.51'
Set-TextValue $ws.Cells.Item(8,7) '0.5919'
$ws.Cells.Item(8,8).Value = 6

# row 9
$aStyleSrc.Copy()
$ws.Cells.Item(9,1).PasteSpecial(-4122)
$ws.Cells.Item(9,1).Value = 7
Set-TextValue $ws.Cells.Item(9,2) '009138'
Set-TextValue $ws.Cells.Item(9,3) '嘉实瑞成两年持有期混合A'
Set-TextValue $ws.Cells.Item(9,4) '22.27'
Set-TextValue $ws.Cells.Item(9,5) '75.95'
Set-TextValue $ws.Cells.Item(9,6) '2.51'
Set-TextValue $ws.Cells.Item(9,7) '0.5590'
$ws.Cells.Item(9,8).Value = 6

# row 10
$aStyleSrc.Copy()
$ws.Cells.Item(10,1).PasteSpecial(-4122)
$ws.Cells.Item(10,1).Value = 8
Set-TextValue $ws.Cells.Item(10,2) '007368'
Set-TextValue $ws.Cells.Item(10,3) '浙商沪港深精选混合A'
Set-TextValue $ws.Cells.Item(10,4) '8.76'
Set-TextValue $ws.Cells.Item(10,5) '92.49'
Set-TextValue $ws.Cells.Item(10,6) '2.87'
Set-TextValue $ws.Cells.Item(10,7) '0.2514'
$ws.Cells.Item(10,8).Value = 10

# row 11
$aStyleSrc.Copy()
$ws.Cells.Item(11,1).PasteSpecial(-4122)
$ws.Cells.Item(11,1).Value = 9
Set-TextValue $ws.Cells.Item(11,2) '900100'
Set-TextValue $ws.Cells.Item(11,3) '中信卓越成长两年持有期混合C'
Set-TextValue $ws.Cells.Item(11,4) '6.91'
Set-TextValue $ws.Cells.Item(11,5) '93.07'
Set-TextValue $ws.Cells.Item(11,6) '3.16'
Set-TextValue $ws.Cells.Item(11,7) '0.2184'
$ws.Cells.Item(11,8).Value = 7

# row 12
$aStyleSrc.Copy()
$ws.Cells.Item(12,1).PasteSpecial(-4122)
$ws.Cells.Item(12,1).Value = 10
Set-TextValue $ws.Cells.Item(12,2) '002653'
Set-TextValue $ws.Cells.Item(12,3) '泰康沪港深精选灵活配置混合'
Set-TextValue $ws.Cells.Item(12,4) '7.89'
Set-TextValue $ws.Cells.Item(12,5) '87.05'
Set-TextValue $ws.Cells.Item(12,6) '2.52'
Set-TextValue $ws.Cells.Item(12,7) '0.1988'
$ws.Cells.Item(12,8).Value = 6

# row 13
$aStyleSrc.Copy()
$ws.Cells.Item(13,1).PasteSpecial(-4122)
$ws.Cells.Item(13,1).Value = 11
Set-TextValue $ws.Cells.Item(13,2) '457001'
Set-TextValue $ws.Cells.Item(13,3) '国富亚洲机会股票 (QDII)'
Set-TextValue $ws.Cells.Item(13,4) '5.93'
Set-TextValue $ws.Cells.Item(13,5) '77.36'
Set-TextValue $ws.Cells.Item(13,6) '2.53'
Set-TextValue $ws.Cells.Item(13,7) '0.1500'
$ws.Cells.Item(13,8).Value = 10

# row 14
$aStyleSrc.Copy()
$ws.Cells.Item(14,1).PasteSpecial(-4122)
$ws.Cells.Item(14,1).Value = 12
Set-TextValue $ws.Cells.Item(14,2) '005335'
Set-TextValue $ws.Cells.Item(14,3) '浙商全景消费混合'
Set-TextValue $ws.Cells.Item(14,4) '2.30'
Set-TextValue $ws.Cells.Item(14,5) '93.36'
Set-TextValue $ws.Cells.Item(14,6) '6.07'
Set-TextValue $ws.Cells.Item(14,7) '0.1396'
$ws.Cells.Item(14,8).Value = 8

# row 15
$aStyleSrc.Copy()
$ws.Cells.Item(15,1).PasteSpecial(-4122)
$ws.Cells.Item(15,1).Value = 13
Set-TextValue $ws.Cells.Item(15,2) '009139'
Set-TextValue $ws.Cells.Item(15,3) '嘉实瑞成两年持有期混合C'
Set-TextValue $ws.Cells.Item(15,4) '4.37'
Set-TextValue $ws.Cells.Item(15,5) '75.95'
Set-TextValue $ws.Cells.Item(15,6) '2.51'
Set-TextValue $ws.Cells.Item(15,7) '0.1097'
$ws.Cells.Item(15,8).Value = 6

# row 16
$aStyleSrc.Copy()
$ws.Cells.Item(16,1).PasteSpecial(-4122)
$ws.Cells.Item(16,1).Value = 14
Set-TextValue $ws.Cells.Item(16,2) '003580'
Set-TextValue $ws.Cells.Item(16,3) '泰康沪港深价值优选灵活配置混合'
Set-TextValue $ws.Cells.Item(16,4) '1.80'
Set-TextValue $ws.Cells.Item(16,5) '85.70'
Set-TextValue $ws.Cells.Item(16,6) '2.48'
Set-TextValue $ws.Cells.Item(16,7) '0.0446'
$ws.Cells.Item(16,8).Value = 7

# row 17
$aStyleSrc.Copy()
$ws.Cells.Item(17,1).PasteSpecial(-4122)
$ws.Cells.Item(17,1).Value = 15
Set-TextValue $ws.Cells.Item(17,2) '012345'
Set-TextValue $ws.Cells.Item(17,3) '嘉实领先优势混合型证券投资基金C'
Set-TextValue $ws.Cells.Item(17,4) '0.82'
Set-TextValue $ws.Cells.Item(17,5) '82.99'
Set-TextValue $ws.Cells.Item(17,6) '3.01'
Set-TextValue $ws.Cells.Item(17,7) '0.0247'
$ws.Cells.Item(17,8).Value = 9

# row 18
$aStyleSrc.Copy()
$ws.Cells.Item(18,1).PasteSpecial(-4122)
$ws.Cells.Item(18,1).Value = 16
Set-TextValue $ws.Cells.Item(18,2) '007369'
Set-TextValue $ws.Cells.Item(18,3) '浙商沪港深精选混合C'
Set-TextValue $ws.Cells.Item(18,4) '0.54'
Set-TextValue $ws.Cells.Item(18,5) '92.49'
Set-TextValue $ws.Cells.Item(18,6) '2.87'
Set-TextValue $ws.Cells.Item(18,7) '0.0155'
$ws.Cells.Item(18,8).Value = 10

# row 19
$aStyleSrc.Copy()
$ws.Cells.Item(19,1).PasteSpecial(-4122)
$ws.Cells.Item(19,1).Value = 17
Set-TextValue $ws.Cells.Item(19,2) '007287'
Set-TextValue $ws.Cells.Item(19,3) '合煦智远消费主题股票A'
Set-TextValue $ws.Cells.Item(19,4) '0.23'
Set-TextValue $ws.Cells.Item(19,5) '83.01'
Set-TextValue $ws.Cells.Item(19,6) '6.58'
Set-TextValue $ws.Cells.Item(19,7) '0.0151'
$ws.Cells.Item(19,8).Value = 3

# row 20
$aStyleSrc.Copy()
$ws.Cells.Item(20,1).PasteSpecial(-4122)
$ws.Cells.Item(20,1).Value = 18
Set-TextValue $ws.Cells.Item(20,2) '010777'
Set-TextValue $ws.Cells.Item(20,3) '浙商智选家居股票A'
Set-TextValue $ws.Cells.Item(20,4) '0.15'
Set-TextValue $ws.Cells.Item(20,5) '90.92'
Set-TextValue $ws.Cells.Item(20,6) '6.83'
Set-TextValue $ws.Cells.Item(20,7) '0.0102'
$ws.Cells.Item(20,8).Value = 6

# row 21
$aStyleSrc.Copy()
$ws.Cells.Item(21,1).PasteSpecial(-4122)
$ws.Cells.Item(21,1).Value = 19
Set-TextValue $ws.Cells.Item(21,2) '007288'
Set-TextValue $ws.Cells.Item(21,3) '合煦智远消费主题股票C'
Set-TextValue $ws.Cells.Item(21,4) '0.05'
Set-TextValue $ws.Cells.Item(21,5) '83.01'
Set-TextValue $ws.Cells.Item(21,6) '6.58'
Set-TextValue $ws.Cells.Item(21,7) '0.0033'
$ws.Cells.Item(21,8).Value = 3

# row 22
$aStyleSrc.Copy()
$ws.Cells.Item(22,1).PasteSpecial(-4122)
$ws.Cells.Item(22,1).Value = 20
Set-TextValue $ws.Cells.Item(22,2) '010778'
Set-TextValue $ws.Cells.Item(22,3) '浙商智选家居股票C'
Set-TextValue $ws.Cells.Item(22,4) '0.03'
Set-TextValue $ws.Cells.Item(22,5) '90.92'
Set-TextValue $ws.Cells.Item(22,6) '6.83'
Set-TextValue $ws.Cells.Item(22,7) '0.0020'
$ws.Cells.Item(22,8).Value = 6

# ============================================================
# Step 2: prepend a "2022-Q1" summary row into "总计" sheet
# ============================================================
$zj = $wb.Worksheets.Item('总计')
$zj.Rows.Item(2).Insert()
$zj.Range("B2:D2").ClearFormats()
$zj.Cells.Item(3,1).Copy()
$zj.Cells.Item(2,1).PasteSpecial(-4122)
$zj.Cells.Item(2,1).Value = 0
$zj.Cells.Item(2,2).Value = '2022-Q1'
$zj.Cells.Item(2,3).Value = 21
$zj.Cells.Item(2,4).Value = 13.99

